$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new Product column values (C2:C4)
$ws.Range("C2").Value = "A"
$ws.Range("C3").Value = "B"
$ws.Range("C4").Value = "C"

# Re-enter B4 as text (quote-prefixed "123") instead of the number 123
$ws.Range("B4").Value = "'123"

# Update the active selection to C4
$ws.Range("C4").Select()
